$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.481.07"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.101.84"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "383.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0852"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "3.594.50"
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "3.104.96"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.992"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "51.518.60"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.168"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0468"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.297"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.71%  "
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "2.059.20"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").Value = "3.413.04"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0326"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
